$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    'D2' = '247.91'
    'D3' = '22.36'
    'D4' = '5.640'
    'D5' = '0.05609'
    'D7' = '6.465'
    'D8' = '0.8018'
    'D9' = '1.062'
    'D10' = '0.1430'
    'D11' = '0.07410'
    'D12' = '0.03197'
    'B13' = 'ProBitToken'
    'C13' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'D13' = '0.1277'
    'E13' = '12ProBitTokenPROB'
    'B14' = 'BitrueCoin'
    'C14' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D14' = '0.02975'
    'E14' = '13BitrueCoinBTR'
    'B15' = 'BitMartToken'
    'C15' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D15' = '0.09263'
    'E15' = '14BitMartTokenBMX'
    'B16' = 'BitForexToken'
    'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D16' = '0.001668'
    'E16' = '15BitForexTokenBF'
    'B17' = 'MCDex'
    'C17' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D17' = '3.252'
    'E17' = '16MCDexMCB'
    'B18' = 'CoinExToken'
    'C18' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'D18' = '0.04681'
    'E18' = '17CoinExTokenCET'
    'B19' = 'One'
    'C19' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'D19' = '0.0005742'
    'E19' = '18OneONEWorstin24h'
    'B20' = 'TigerCash'
    'C20' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D20' = '0.006270'
    'E20' = '19TigerCashTCH'
    'B21' = 'BitKan'
    'C21' = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
    'D21' = '0.001053'
    'E21' = '20BitKanKAN'
    'B22' = 'HotbitToken'
    'C22' = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
    'D22' = '0.003811'
    'E22' = '21HotbitTokenHTB'
    'B23' = 'NitroEx'
    'C23' = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
    'D23' = '0.0001500'
    'E23' = '22NitroExNTX'
    'B24' = 'UpBots'
    'C24' = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
    'D24' = '0.0004602'
    'E24' = '23UpBotsUBXT'
    'B25' = 'LEO'
    'C25' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D25' = '3.979'
    'E25' = '24LEOLEO'
    'B26' = 'BTSEToken'
    'C26' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D26' = '2.116'
    'E26' = '25BTSETokenBTSE'
    'B27' = 'BitpandaEcosystemToken'
    'C27' = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
    'D27' = '0.3311'
    'E27' = '26BitpandaEcosystemTokenBEST'
    'D40' = '0.04192'
    'D41' = '0.007031'
    'E41' = '40KickTokenKICK'
    'D44' = '0.008805'
    'D45' = '0.00005670'
    'D47' = '0.6803'
    'D48' = '0.02929'
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = "'" + $changes[$cell]
    $ws.Range($cell).Style = "Normal"
}

Write-Host "Applied $($changes.Count) cell changes"
